# Update res_bus vm_pu values for rows 2-25 (data rows 0-23) per the new
# 380 kV case results. Only columns B:F and I:N change; column G stays 1
# and column H (empty) / column A (index) are untouched.

$ws = $excel.ActiveWorkbook.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

$data = @"
2:1.02,1.02097251208699,1.025227522138852,1.046499714841461,1.050236659084414,1.027652599957388,1.026166605366576,1.028053709786372,1.049264978267421,1.052991496207353,1.027623879649152
3:1.02,1.021904087663062,1.025873153972127,1.047771839125077,1.051580900851814,1.027754859853377,1.026735692402735,1.028507039143563,1.05034756155312,1.054146773395805,1.028193774854177
4:1.02,1.022506935067554,1.02629074800314,1.048595713681215,1.052451500481056,1.027819572241009,1.02710339399349,1.028799516947302,1.051048236949373,1.0548945607812,1.028561998623367
5:1.02,1.022760385045016,1.026466262032566,1.048942245750696,1.052817690229949,1.027846428196417,1.027257846994057,1.028922268848545,1.051342843355426,1.055208990787999,1.02871667096491
6:1.02,1.022802941157726,1.026495729086048,1.049000440331395,1.052879186270475,1.027850916944026,1.02728377276857,1.028942867346341,1.051392311598102,1.055261788538541,1.028742633556997
7:1.02,1.02251032162827,1.026293093398588,1.048600343369208,1.052456392779064,1.027819932463977,1.027105458308548,1.028801157973621,1.051052173325975,1.054898761972143,1.028564065869989
8:1.02,1.021287330632618,1.025445751588003,1.046929486040021,1.050690790273903,1.027687460515932,1.02635904164137,1.028207091513833,1.049630807653294,1.053381877457932,1.027816589205529
9:1.02,1.01913270860595,1.023951350769963,1.043990702732765,1.047585513527534,1.027442890548744,1.025039677671849,1.027153740327915,1.047127432267162,1.050710755339201,1.026495351587837
10:1.02,1.017696614033007,1.02295428922285,1.042035095711526,1.045519217413411,1.027272379955651,1.024157383540399,1.02644715541768,1.045459269277532,1.04893114619275,1.02561180449752
11:1.02,1.017074848851538,1.022522375217032,1.041189122321459,1.044625386519218,1.027196781571935,1.023774700821082,1.026140172382716,1.044737094937413,1.048160806118585,1.025228578324847
12:1.02,1.016843908784027,1.022361917121736,1.040875010832447,1.044293509321804,1.027168435979471,1.023632459077601,1.026025991491831,1.044468868140997,1.047874702147725,1.025086134581776
13:1.02,1.016893445702186,1.022396337088432,1.040942383385163,1.044364692196038,1.027174528190704,1.02366297476867,1.026050490627421,1.044526402782303,1.047936070882224,1.025116693608626
14:1.02,1.017055759043775,1.022509112225463,1.041163155317802,1.044597950778633,1.027194443919986,1.023762945037583,1.026130737287177,1.044714922810705,1.048137155978963,1.025216805846787
15:1.02,1.017155767141319,1.02257859329386,1.041299196122291,1.044741686436317,1.027206679555557,1.023824527263557,1.026180159550426,1.044831078938492,1.048261055679294,1.025278475526586
16:1.02,1.017737880100698,1.022982950231333,1.042091257220419,1.04557855657088,1.027277359990993,1.024182767402757,1.026467507280992,1.045507200641655,1.048982276081267,1.025637224407875
17:1.02,1.018103043932715,1.023236545236681,1.042588313565403,1.046103739569858,1.027321223361358,1.024407309926148,1.026647478351398,1.045931353134016,1.049434742226267,1.025862085807404
18:1.02,1.018316044871541,1.023384445499806,1.042878317125771,1.046410155844231,1.027346637689091,1.024538219666188,1.026752353291813,1.04617876882096,1.049698681584422,1.02599318145429
19:1.02,1.018388673877557,1.023434872683182,1.042977214309937,1.046514650513997,1.027355274398519,1.024582845958336,1.026788096102475,1.046263133727433,1.049788682127104,1.026037871120892
20:1.02,1.018063864578264,1.023209338686385,1.042534975966877,1.046047383532951,1.027316534863008,1.02438322506073,1.026628179408992,1.045885844076639,1.049386194471859,1.025837966738714
21:1.02,1.017007961499308,1.022475903468195,1.041098140171342,1.044529258330091,1.027188586546138,1.023733508941874,1.026107110884785,1.044659407798463,1.048077940488237,1.025187327948447
22:1.02,1.016344138252725,1.022014613687925,1.040195441743549,1.043575511900692,1.027106607407403,1.023324449389974,1.025778605315533,1.04388841870618,1.047255588770356,1.024777687485062
23:1.02,1.016696037294373,1.022259166007111,1.040673913868998,1.044081039692024,1.02715021128608,1.023541352304588,1.025952836415031,1.044297123835323,1.047691514663763,1.024994898426695
24:1.02,1.01808156802019,1.023221632209236,1.0425590766871,1.04607284813417,1.027318653920225,1.024394108166571,1.026636900076682,1.045906407614799,1.049408131039222,1.025848865299814
25:1.02,1.019689675183869,1.024337833856198,1.044749809791122,1.048387608917724,1.027507435091734,1.025381245917533,1.027426826610438,1.047774475473954,1.051401097383402,1.026837404899609
"@

foreach ($line in ($data -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ":"
    $row = [int]$parts[0]
    $values = $parts[1] -split ","
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = [double]$values[$i]
        $ws.Cells.Item($row, $col).Value = $val
    }
}
